# Update version string references from the Jan 30 build to the
# "version 1.0.0 (Feb 3 2026)" release, rebuilt Feb 03 2026 10.14.00 EST.

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: ..."
$wsAbout.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended Citation text embeds the version string.
$oldCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Kostromovskaya Coal Mine, Russia, M1343, version '" + $oldVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"
$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Kostromovskaya Coal Mine, Russia, M1343, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"
$wsAbout.Range("A6").Value = $newCitation

# Column S (rows 2-14) on the data sheet holds the same version string.
for ($r = 2; $r -le 14; $r++) {
    $cell = $wsData.Cells.Item($r, 19)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
